# Add an "OrcID" row to the author profile form, per commit:
# "...author details form to include OrcIDs..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 10 (Bio row), shifting everything below
# down by one. Excel automatically re-targets hyperlinks, merged cells and
# the used-range dimension when rows are inserted this way.
$ws.Rows("10:10").Insert()

# The freshly inserted row has no formatting of its own; pull the label-cell
# style (thin/medium border combo, no fill) from the row directly below
# (which used to be row 10, the "Bio" row) so the new "OrcID" label matches
# the other label cells in this section.
$ws.Range("A11").Copy()
$ws.Range("A10").PasteSpecial(-4122)

# Populate the new row.
$ws.Range("A10").Value = "OrcID"
$ws.Range("B10").Value = "0000-0003-1105-4971"

# Give the OrcID value a distinct, larger bold font (Arial 14, black) to make
# it stand out, matching the rest of the form's emphasis styling.
$valueFont = $ws.Range("B10").Font
$valueFont.Name = "Arial"
$valueFont.Size = 14
$valueFont.Bold = $true
$valueFont.Color = 0

# Give the row a bit more breathing room for the larger font.
$ws.Rows("10:10").RowHeight = 17.4

# Leave the cursor parked on the new value cell, as the author did.
[void]$ws.Range("B10").Select()
